# "Generate Report for Archive"
#
# 1) The per-language localization status text changes from
#    "Ready for handoff" to "In Translation" everywhere it appears:
#      - Overview sheet: E2, F2, E3, F3 (zh-cn / de-de status columns)
#      - zh-cn sheet:     C2, C3 (Status column)
#      - de-de sheet:     C2, C3 (Status column)
#
# 2) The (now shorter) status columns are narrowed to match the new
#    content width:
#      - Overview sheet: columns E and F
#      - zh-cn sheet:     column C
#      - de-de sheet:     column C

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Update status values ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Narrow the status columns to their new content width ---
$newColumnWidth = 13.4101848602295

$wsOverview.Range("E1").EntireColumn.ColumnWidth = $newColumnWidth
$wsOverview.Range("F1").EntireColumn.ColumnWidth = $newColumnWidth
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = $newColumnWidth
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = $newColumnWidth
